$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to include units
$ws.Range("B1").Value = "Python (ms)"
$ws.Range("C1").Value = "C (ms)"

# Update Python timing values (converted from ms-rounded ints to precise ms values)
$ws.Range("B2").Value = 3.66127
$ws.Range("B3").Value = 3.0455299999999998
$ws.Range("B4").Value = 3.3317299999999999
$ws.Range("B5").Value = 2.4952700000000001
$ws.Range("B6").Value = 3.6767500000000002

# Widen column B to fit the new, longer header text
$ws.Columns("B").ColumnWidth = 11.67

# Move the active selection to C1
$ws.Range("C1").Select()
